$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values (F14/N14 and F15/N15)
$ws.Range("F14").Value = 82
$ws.Range("N14").Value = 85
$ws.Range("F15").Value = 80
$ws.Range("N15").Value = 85

# Move the active selection from G14 to A14, matching the saved sheet view
$ws.Range("A14").Select()
